$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/disposition-reason"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- Elements sheet update ---
$elements = $wb.Worksheets.Item("Elements")

# Row 2 is the root "Extension" element; its Constraint(s) column (AI) had
# the ele-1/ext-1 constraint text duplicated from the child Extension.extension
# row. That duplication is removed here.
$elements.Range("AI2").Value = ""

# Extension.url's "Fixed Value" (Q5) shares the same underlying text as the
# workbook's URL (Metadata!B2) -- both cells held the literal
# "http://ibm.com/..." URL string, so both must be updated together to stay
# in sync with the new "http://linuxforhealth.org/..." URL.
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/disposition-reason"
